# feat: add 2022-Q1 data
#
# Before: sheets are  2021-Q3 | 2021-Q4 | 总计
# After:  sheets are  2021-Q3 | 2021-Q4 | 2022-Q1 | 总计
#
# The existing "总计" sheet (which held the latest-quarter fund detail,
# i.e. the same shape/style as "2021-Q4") is renamed to "2022-Q1" and
# its data is replaced with the new quarter's fund holdings. A brand
# new "总计" sheet is then inserted right after it, recreating the
# running summary table (日期 / 持有数量(只) / 持有市值(亿元)) with an
# extra leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# "2021-Q4" already has exactly the layout/styling we need to
# reproduce for "2022-Q1" (and for the header/index-column look of
# the summary table), so use it as the formatting template throughout.
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------
# 1. Rename "总计" -> "2022-Q1" and overwrite its content with the
#    new quarter's fund-holding detail.
# ---------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q1.Range($cols[$i] + "1")
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
}
$q1.Range("B1:H1").Style = "Normal"
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$template.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$fundRow = @("516910", "富国中证现代物流交易型开放式指数证券投资基金", "0.42", "98.43", "4.53", "0.0190")
for ($i = 0; $i -lt $fundRow.Length; $i++) {
    $cell = $q1.Range($cols[$i] + "2")
    $cell.NumberFormat = "@"
    $cell.Value = $fundRow[$i]
    $cell.Style = "Normal"
}
$q1.Range("H2").Value = 5

# ---------------------------------------------------------------
# 2. Insert a fresh "总计" sheet right after "2022-Q1" and rebuild
#    the running summary table (new quarter on top, older quarters
#    shifted down).
# ---------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$totalCols = @("B", "C", "D")
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $cell = $total.Range($totalCols[$i] + "1")
    $cell.NumberFormat = "@"
    $cell.Value = $totalHeaders[$i]
}
$total.Range("B1:D1").Style = "Normal"
$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q1", 1, 0.02),
    @(1, "2021-Q4", 1, 0.02),
    @(2, "2021-Q3", 1, 0.03)
)
for ($r = 0; $r -lt $summaryRows.Length; $r++) {
    $rowNum = $r + 2
    $total.Range("A$rowNum").Value = $summaryRows[$r][0]
    $template.Range("A2").Copy()
    $total.Range("A$rowNum").PasteSpecial(-4122)
    $total.Range("B$rowNum").Value = $summaryRows[$r][1]
    $total.Range("C$rowNum").Value = $summaryRows[$r][2]
    $total.Range("D$rowNum").Value = $summaryRows[$r][3]
}
